# Updates the cryptocurrency price/volume table (and a couple of re-ranked rows)
# to match the latest scrape, as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.987.93'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '3.858.74'
$ws.Range("E3").Value = '  +1.41%  '

$ws.Range("D5").Value = '''697.48'
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").Value = '''173.47'
$ws.Range("E6").Value = '  +0.49%  '

$ws.Range("D7").Value = '3.856.07'
$ws.Range("E7").Value = '  +1.38%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = '''0.525'
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").Value = '''0.162'
$ws.Range("E10").Value = '  -0.17%  '

$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("D12").Value = '''0.459'
$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '''0.0000258'
$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("D14").Value = '''36.33'
$ws.Range("E14").Value = '  +0.07%  '

$ws.Range("D15").Value = '4.511.25'
$ws.Range("E15").Value = '  +1.45%  '

$ws.Range("D16").Value = '3.859.52'
$ws.Range("E16").Value = '  +1.25%  '

$ws.Range("D17").Value = '71.050.30'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").Value = '''17.59'
$ws.Range("E18").Value = '  -0.80%  '

$ws.Range("D19").Value = '''7.20'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("D21").Value = '''10.82'
$ws.Range("E21").Value = '  -4.45%  '

$ws.Range("D22").Value = '''498.56'
$ws.Range("E22").Value = '  +3.95%  '

$ws.Range("D23").Value = '''0.720'
$ws.Range("E23").Value = '  +0.56%  '

$ws.Range("E24").Value = '  +4.08%  '

$ws.Range("E25").Value = '  +1.32%  '

$ws.Range("D26").Value = '''10.66'
$ws.Range("E26").Value = '  +2.18%  '

$ws.Range("D27").Value = '''12.19'
$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").Value = '''2.13'
$ws.Range("E28").Value = '  -1.74%  '

$ws.Range("E29").Value = '  +1.72%  '

$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("D31").Value = '''7.55'
$ws.Range("E31").Value = '  +0.53%  '

$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '''0.183'
$ws.Range("E33").Value = '  +3.72%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''29.57'
$ws.Range("E34").Value = '  +0.15%  '

$ws.Range("D35").Value = '''9.20'
$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("D36").Value = '3.813.59'
$ws.Range("E36").Value = '  +1.54%  '

$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("E38").Value = '  +1.73%  '

$ws.Range("D39").Value = '''2.37'
$ws.Range("E39").Value = '  +6.73%  '

$ws.Range("E40").Value = '  +8.34%  '

$ws.Range("E41").Value = '  -2.08%  '

$ws.Range("D42").Value = '''6.02'
$ws.Range("E42").Value = '  +0.58%  '

$ws.Range("E43").Value = '  -0.04%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = '''0.000313'
$ws.Range("E45").Value = '  -6.55%  '

$ws.Range("D46").Value = '''163.76'
$ws.Range("E46").Value = '  +1.91%  '

$ws.Range("D47").Value = '''49.22'
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").Value = '''417.46'
$ws.Range("E48").Value = '  +3.31%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '''1.39'
$ws.Range("E49").Value = '  -1.25%  '

$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = '''0.301'
$ws.Range("E50").Value = '  +0.81%  '

$ws.Range("D51").Value = '''43.50'
$ws.Range("E51").Value = '  -4.36%  '
